{"js": "// Regenerate the lattice-multiplication exercise grid with a new set of\n// problems (same 5x3 table shape; each cell's title / top digits / rule /\n// left digits get replaced with a new multiplication problem).\nconst target = [\n  [\"21 x 35\\u000b  3    5\\u000b  ----\\u000b2|    |\\u000b1|    |\", \"93 x 92\\u000b  9    2\\u000b  ----\\u000b9|    |\\u000b3|    |\", \"22 x 62\\u000b  6    2\\u000b  ----\\u000b2|    |\\u000b2|    |\"],\n  [\"20 x 61\\u000b  6    1\\u000b  ----\\u000b2|    |\\u000b0|    |\", \"45 x 69\\u000b  6    9\\u000b  ----\\u000b4|    |\\u000b5|    |\", \"99 x 19\\u000b  1    9\\u000b  ----\\u000b9|    |\\u000b9|    |\"],\n  [\"13 x 40\\u000b  4    0\\u000b  ----\\u000b1|    |\\u000b3|    |\", \"28 x 10\\u000b  1    0\\u000b  ----\\u000b2|    |\\u000b8|    |\", \"88 x 85\\u000b  8    5\\u000b  ----\\u000b8|    |\\u000b8|    |\"],\n  [\"48 x 10\\u000b  1    0\\u000b  ----\\u000b4|    |\\u000b8|    |\", \"87 x 86\\u000b  8    6\\u000b  ----\\u000b8|    |\\u000b7|    |\", \"87 x 21\\u000b  2    1\\u000b  ----\\u000b8|    |\\u000b7|    |\"],\n  [\"49 x 99\\u000b  9    9\\u000b  ----\\u000b4|    |\\u000b9|    |\", \"80 x 51\\u000b  5    1\\u000b  ----\\u000b8|    |\\u000b0|    |\", \"96 x 15\\u000b  1    5\\u000b  ----\\u000b9|    |\\u000b6|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\n// Replace the text of the existing paragraph's range in place (rather than\n// clearing + re-inserting the whole cell body) so the run-level formatting\n// that is already on the text (sz=32) is preserved exactly.\nfor (let r = 0; r < target.length && r < table.rowCount; r++) {\n  for (let c = 0; c < target[r].length && c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    const para = cell.body.paragraphs.items[0];\n    const range = para.getRange();\n    range.insertText(target[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Regenerate the lattice-multiplication exercise grid with a new set of\n# problems (same 5x3 table shape; each cell's title / top-digit row /\n# separator / two left-digit rows get replaced with a new multiplication\n# problem). The table shape (5 rows x 3 columns) and run formatting\n# (sz=32) are left untouched; only the cell text changes.\n#\n# `v below is PowerShell's escape sequence for a vertical-tab character,\n# which Word's Range.Text setter turns into a manual line break (<w:br/>)\n# matching the existing ``<w:t>...<w:br/>`` run layout in each cell.\n$target = @(\n    @(\"21 x 35`v  3    5`v  ----`v2|    |`v1|    |\", \"93 x 92`v  9    2`v  ----`v9|    |`v3|    |\", \"22 x 62`v  6    2`v  ----`v2|    |`v2|    |\"),\n    @(\"20 x 61`v  6    1`v  ----`v2|    |`v0|    |\", \"45 x 69`v  6    9`v  ----`v4|    |`v5|    |\", \"99 x 19`v  1    9`v  ----`v9|    |`v9|    |\"),\n    @(\"13 x 40`v  4    0`v  ----`v1|    |`v3|    |\", \"28 x 10`v  1    0`v  ----`v2|    |`v8|    |\", \"88 x 85`v  8    5`v  ----`v8|    |`v8|    |\"),\n    @(\"48 x 10`v  1    0`v  ----`v4|    |`v8|    |\", \"87 x 86`v  8    6`v  ----`v8|    |`v7|    |\", \"87 x 21`v  2    1`v  ----`v8|    |`v7|    |\"),\n    @(\"49 x 99`v  9    9`v  ----`v4|    |`v9|    |\", \"80 x 51`v  5    1`v  ----`v8|    |`v0|    |\", \"96 x 15`v  1    5`v  ----`v9|    |`v6|    |\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 0; $r -lt $target.Length -and $r -lt $rowCount; $r++) {\n    $rowValues = $target[$r]\n    for ($c = 0; $c -lt $rowValues.Length -and $c -lt $colCount; $c++) {\n        $cell = $t.Cell($r + 1, $c + 1)\n        # Assigning Range.Text in place preserves the existing run's\n        # formatting (sz=32) instead of resetting to document defaults.\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
